# "Ajout de l'importation des données du fichier xls et aperçu"
#
# A new "preview" column is inserted right after column A (NOM): the
# imported/raw values that used to live in column A are duplicated into
# the new column B, the old column B (Prénom) shifts right to column C,
# and the new header cell A1 is relabelled "test" (the import preview
# marker), while the old header "NOM" slides into B1 and "Prénom" into C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column A (values + formatting) and insert it before column B —
# this is Excel's "Insert Copied Cells" gesture: it shifts the existing
# column B (and everything right of it) over to column C, and the newly
# inserted column B is a duplicate of column A.
[void]$ws.Columns("A").Copy()
[void]$ws.Columns("B").Insert()

# Re-label the header of the new first column.
$ws.Range("A1").Value2 = "test"

# Leave the selection where the user ended up after the import.
[void]$ws.Range("A2").Select()
